$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Actual Consumption (MW)" values for rows 2..193 (data rows for
# 21.01.2026 and 22.01.2026, which become 23.01.2026 and 24.01.2026).
$newB = @(6484,6483,6456,6402,6339,6325,6271,6230,6220,6203,6228,6237,6270,6312,6320,6376,6512,6536,6643,6634,6849,6989,7200,7406,7780,7961,8119,8194,8428,8553,8725,8841,8871,8868,8879,8885,8803,8826,8770,8689,8588,8576,8596,8519,8417,8454,8400,8442,8655,8630,8588,8602,8577,8526,8505,8467,8252,8265,8237,8296,8350,8333,8353,8403,8525,8584,8732,8760,8748,8710,8786,8686,8587,8608,0,8555,8455,8464,8352,8252,8063,8017,7915,7782,7551,7428,7278,7198,7037,6925,6791,6760,6566,6556,6434,6396,6320,6219,6208,6156,6211,6160,6151,6123,6148,6111,6086,6106,6133,0,6108,6139,6164,6166,6200,6242,6259,6320,6329,6411,6536,6619,6697,6805,6907,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)

$firstRow = 2
$lastRow = 193

for ($row = $firstRow; $row -le $lastRow; $row++) {
    $i = $row - $firstRow

    $oldA = $ws.Cells.Item($row, 1).Value2
    $oldD = $ws.Cells.Item($row, 4).Value2

    $newA = $oldA + 2
    $newD = $oldD.Replace("21.01.2026", "23.01.2026").Replace("22.01.2026", "24.01.2026")

    $ws.Cells.Item($row, 1).Value = $newA
    $ws.Cells.Item($row, 2).Value = $newB[$i]
    $ws.Cells.Item($row, 4).Value = $newD
}
